# The SQL-query section of the write-up was reworked; the old lead-in
# page break plus the four paragraphs describing "khóa học / kỳ học /
# sinh viên / niên khóa" are dropped from the document entirely. The
# empty paragraph that used to follow them is left in place.

$d = $word.ActiveDocument
$count = $d.Paragraphs.Count

$startText = "Mỗi năm là một khóa học"
$endText = "Mỗi sinh viên sẽ có 5 môn học tự chọn trong chương trình đào tạo."

$findStart = $d.Content.Duplicate
$foundStart = $findStart.Find.Execute($startText, $false, $false, $false, $false, `
                         $false, $true, 1, $false, "", 0)
$startPos = $findStart.Start

$findEnd = $d.Content.Duplicate
$foundEnd = $findEnd.Find.Execute($endText, $false, $false, $false, $false, `
                       $false, $true, 1, $false, "", 0)
$endPos = $findEnd.Start

$startIdx = -1
$endIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $pStart = $p.Range.Start
    $pEnd = $p.Range.End
    if (($startPos -ge $pStart) -and ($startPos -lt $pEnd)) {
        $startIdx = $i
    }
    if (($endPos -ge $pStart) -and ($endPos -lt $pEnd)) {
        $endIdx = $i
    }
}

# Also swallow the lone page-break paragraph immediately preceding the
# "Mỗi năm là một khóa học" paragraph.
$pageBreakIdx = $startIdx - 1

$deleteStart = $d.Paragraphs.Item($pageBreakIdx).Range.Start
$deleteEnd = $d.Paragraphs.Item($endIdx).Range.End

$rng = $d.Range($deleteStart, $deleteEnd)
$rng.Delete()
